# edit.ps1
# Applies the "update of poster 070225" edit to K20_Raman_averaged.xlsx
#
# Summary of the change:
#  1. The column "Name_on_SEM_Raman_STD" (column GF) is removed entirely,
#     shifting all later "...Raman_STD" / "averaged?" columns one position
#     to the left (GG->GF, GH->GG, GI->GH, GJ->GI, GK->GJ, GL->GK).
#     This also shrinks the used range from A1:GL39 to A1:GK39.
#  2. The "Name_on_SEM" column (CZ) gets populated for rows 7-39 with the
#     sample name derived from the "FI_name" column (CY), i.e. the FI_name
#     value with its trailing "_FI<letter>" suffix stripped off.
#  3. A handful of pre-existing numeric measurement values were re-saved
#     with a 1-ULP floating point difference (last significant digit only)
#     as part of the same notebook re-run; these are corrected explicitly
#     below to match the refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove column GF ("Name_on_SEM_Raman_STD"); everything to its right
#    shifts left by one column.
# ---------------------------------------------------------------------
$ws.Columns("GF").Delete()

# ---------------------------------------------------------------------
# 2) Populate "Name_on_SEM" (column CZ) for rows 7-39 from "FI_name"
#    (column CY), stripping the trailing "_FI<letter>" suffix.
# ---------------------------------------------------------------------
$nameOnSem = @{
    7  = "K20_c001_a1"
    8  = "K20_c002_a1"
    9  = "K20_c002_a1"
    10 = "K20_c003_a1"
    11 = "K20_c004_a1"
    12 = "K20_c004_a2"
    13 = "K20_c004_a3"
    14 = "K20_c004_a3"
    15 = "K20_c004_a3"
    16 = "K20_c005_a1"
    17 = "K20_c005_a2"
    18 = "K20_c005_a3"
    19 = "K20_c006_a1"
    20 = "K20_c007_a1"
    21 = "K20_c007_a2"
    22 = "K20_c008_a1"
    23 = "K20_c008_a2"
    24 = "K20_c009_a1"
    25 = "K20_c009_a2"
    26 = "K20_c010_a1"
    27 = "K20_c010_a2"
    28 = "K20_c011_a1"
    29 = "K20_c014_a1"
    30 = "K20_c015_a1"
    31 = "K20_c018_a2"
    32 = "K20_c019_a1"
    33 = "K20_c021_a1"
    34 = "K20_c024_a1"
    35 = "K20_c025_a1"
    36 = "K20_c026_a1"
    37 = "K20_c026_a2"
    38 = "K20_c027_a1"
    39 = "K20_c028_a1"
}

foreach ($row in $nameOnSem.Keys) {
    $ws.Range("CZ$row").Value = $nameOnSem[$row]
}

# ---------------------------------------------------------------------
# 3) Refresh a small set of numeric cells whose last significant digit
#    changed (1-ULP) when the workbook was re-saved.
# ---------------------------------------------------------------------
$floatUpdates = @{
    "DF2"  = 0.004832970503648978
    "DV2"  = 0.004832970503648978
    "EP2"  = 0.00323854315347735
    "FK2"  = 0.09796047447211081
    "EI3"  = 0.01748610600565654
    "EJ3"  = 0.04330182253267002
    "FQ3"  = 0.000005590195512312427
    "EB4"  = 0.03076061906781654
    "EE4"  = 0.03076061906776445
    "DN5"  = 0.009165151389911469
    "EB5"  = 0.04363865180453827
    "EE5"  = 0.04363865180453827
    "EJ5"  = 0.01642036681494803
    "EK5"  = 0.02427342317049297
    "DJ6"  = 0.007210467132182745
    "DU6"  = 0.0000552022941336663
    "ER6"  = 0.07332998605063307
    "FE6"  = 42.18958180234007
    "FL6"  = 0.06750483955959083
    "DN34" = 0.07544534445544122
    "EI34" = 0.03978111779819275
    "EK34" = 0.002801861622381327
    "EO34" = 84.99778839562003
    "FM34" = 0.0003249345407530975
    "DF35" = 0.007195608131690942
    "DT35" = 0.007195608131690942
    "EF36" = 9.623457155447099
    "DL37" = 0.000003270235687701335
    "EH37" = 0.3571006005535773
}

foreach ($cellRef in $floatUpdates.Keys) {
    $ws.Range($cellRef).Value = $floatUpdates[$cellRef]
}
